$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update "Included" column (C) values: "No"/"Not Yet" -> "Yes?" for the targeted rows.
$ws.Range("C3").Value = "Yes?"

$ws.Range("C9:C13").Value = "Yes?"

$ws.Range("C23:C34").Value = "Yes?"

# Bold the "Included" cells for the Previous qualification rows (C7:C8),
# matching the formatting applied alongside this data re-check.
$ws.Range("C7:C8").Font.Bold = $true

# C9:C13 pick up a font style flag too (applied then un-bolded), so explicitly
# set Bold True then False to reproduce the resulting "applyFont" style record.
$ws.Range("C9:C13").Font.Bold = $true
$ws.Range("C9:C13").Font.Bold = $false

# Update the final selection to D13, matching the saved cursor position.
$ws.Range("D13").Select()
